# Edits to Intro chapter
$d = $word.ActiveDocument
$rsquo = [char]0x2019

# --- 1. First-line indent on the opening paragraph ("My thesis project was designed...") ---
$d.Paragraphs(2).Range.ParagraphFormat.FirstLineIndent = 36

# --- 2. "My project was designed" -> "My thesis project was designed" ---
$d.Paragraphs(2).Range.Find.Execute("My project was designed", $true, $false, $false, $false, $false, `
    $true, 1, $false, "My thesis project was designed", 2) | Out-Null

# --- 3. Drop the stray "2006" before "Minnesota State Wildlife Action Plan" ---
$d.Paragraphs(2).Range.Find.Execute("Minnesota Department of Natural Resources (hereafter, MN NDR) 2006 Minnesota State Wildlife Action Plan", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Minnesota Department of Natural Resources (hereafter, MN NDR) Minnesota State Wildlife Action Plan", 2) | Out-Null

# --- 4. "in which it was identified" -> "where the Anoka Sand Plain was identified" ---
$d.Paragraphs(2).Range.Find.Execute("Anoka Sand Plain subsection profile, in which it was identified as containing", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Anoka Sand Plain subsection profile, where the Anoka Sand Plain was identified as containing", 2) | Out-Null

# --- 5. "39 of which are" -> "39 of these species are" ---
$d.Paragraphs(2).Range.Find.Execute("within the subsection; 39 of which are federal or state endangered", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "within the subsection; 39 of these species are federal or state endangered", 2) | Out-Null

# --- 6/7. Move "since the 1930s" earlier, drop parens around pine species, "though" -> "although" ---
$d.Paragraphs(3).Range.Find.Execute("timber production and recreation, and much of the area", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "timber production and recreation since the 1930s, and much of the area", 2) | Out-Null

$search1930 = "Pinus strobus and Pinus resinosa) since the 1930" + $rsquo + "s with the original intent"
$d.Paragraphs(3).Range.Find.Execute($search1930, `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Pinus strobus and Pinus resinosa) with the original intent", 2) | Out-Null

$d.Paragraphs(3).Range.Find.Execute("soil during periods of drought, though timber profitability", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "soil during periods of drought, although timber profitability", 2) | Out-Null

# --- 8. Drop parentheses around "such as oaks" ---
$d.Paragraphs(5).Range.Find.Execute("hardwood species (such as oaks) had been removed", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "hardwood species such as oaks had been removed", 2) | Out-Null

# --- 9. "the scope of the study" -> "the original scope of my study" ---
$d.Paragraphs(5).Range.Find.Execute("significantly compromised the scope of the study to directly", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "significantly compromised the original scope of my study to directly", 2) | Out-Null

# --- 10. Widen left margin from 1440 -> 2160 twips (72pt -> 108pt) ---
$d.PageSetup.LeftMargin = 108

# --- 11. New closing paragraph about Chapter 2 / collaboration ---
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newParaText = "Chapter 2 of my thesis was written to conform to the submission requirements of the " + `
    "Journal of Ecological Restoration" + `
    ", where I intend to submit my thesis. I wrote Chapter 2 in collaboration with my academic advisor Dr. Todd Arnold and Dr. Althea " + `
    "ArchMiller" + `
    ", both of whom will be co-authors on any resulting publication. My use of plural pronouns throughout Chapter 2 is reflective of this collaboration, but I take individual responsibility for the full content of my thesis."
$newPara.Range.InsertAfter($newParaText)

# Italicize the journal title within the newly added paragraph.
$italicRange = $d.Paragraphs.Last.Range
$italicRange.Find.Execute("Journal of Ecological Restoration") | Out-Null
$italicRange.Font.Italic = 1

# Recreate the "_GoBack" bookmark Word leaves at the last edit point (best effort;
# a no-op if bookmark creation isn't supported by this host).
try {
    $goBackRange = $d.Paragraphs.Last.Range
    $goBackRange.Find.Execute(", where I intend to submit my thesis.") | Out-Null
    $d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
} catch {}

Write-Output "text edits applied"
